# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment schedule"
# sheet, shifting the old "Late"/"heading"/"Outstanding" columns one slot
# to the right, then make "Repayment schedule" the active sheet/tab with
# the cursor on K18 (previously the "NewLoanInput" sheet/tab was active).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column at position 14 (N); existing N/O/P data shifts to O/P/Q.
$ws.Columns.Item(14).Insert() | Out-Null

# Match the column width Excel shows for the freshly inserted column.
$ws.Columns.Item(14).ColumnWidth = 10.14

# Make "Repayment schedule" the active sheet/tab and move the selection.
$ws.Activate() | Out-Null
$ws.Range("K18").Select() | Out-Null
